# Updates cryptos list: refreshed Price (D) and Volume(1h) (E) columns,
# and re-ranked rows 26/27 (InternetComputer(DFINITY) now ranks above WrappedeETH).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.697.96"
$ws.Range("E2").Value = "  +4.15%  "
$ws.Range("D3").Value = "3.628.35"
$ws.Range("E3").Value = "  +2.75%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'630.27"
$ws.Range("E5").Value = "  +3.35%  "
$ws.Range("D6").Value = "'160.15"
$ws.Range("E6").Value = "  +5.02%  "
$ws.Range("D7").Value = "3.626.90"
$ws.Range("E7").Value = "  +2.76%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "'0.495"
$ws.Range("E9").Value = "  +2.59%  "
$ws.Range("D10").Value = "'0.149"
$ws.Range("E10").Value = "  +6.76%  "
$ws.Range("E11").Value = "  +5.49%  "
$ws.Range("D12").Value = "'0.441"
$ws.Range("E12").Value = "  +3.60%  "
$ws.Range("D13").Value = "'0.0000229"
$ws.Range("E13").Value = "  +4.24%  "
$ws.Range("D14").Value = "'33.42"
$ws.Range("E14").Value = "  +5.61%  "
$ws.Range("D15").Value = "4.236.67"
$ws.Range("E15").Value = "  +2.57%  "
$ws.Range("D16").Value = "3.625.81"
$ws.Range("E16").Value = "  +2.92%  "
$ws.Range("D17").Value = "69.825.78"
$ws.Range("E17").Value = "  +4.29%  "
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("D19").Value = "'6.65"
$ws.Range("E19").Value = "  +5.85%  "
$ws.Range("D20").Value = "'16.03"
$ws.Range("E20").Value = "  +4.40%  "
$ws.Range("D21").Value = "'10.22"
$ws.Range("E21").Value = "  +11.11%  "
$ws.Range("D22").Value = "'463.06"
$ws.Range("E22").Value = "  +4.31%  "
$ws.Range("D23").Value = "'0.644"
$ws.Range("E23").Value = "  +2.41%  "
$ws.Range("D24").Value = "'78.62"
$ws.Range("E24").Value = "  +1.43%  "
$ws.Range("D25").Value = "'0.0000137"
$ws.Range("E25").Value = "  +11.27%  "
$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").Value = "'10.76"
$ws.Range("E26").Value = "  +5.78%  "
$ws.Range("B27").Value = "WrappedeETH"
$ws.Range("C27").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D27").Value = "3.766.84"
$ws.Range("E27").Value = "  +2.63%  "
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("D29").Value = "'9.27"
$ws.Range("E29").Value = "  +13.33%  "
$ws.Range("D30").Value = "'2.64"
$ws.Range("E30").Value = "  +4.30%  "
$ws.Range("D31").Value = "'1.72"
$ws.Range("E31").Value = "  +3.29%  "
$ws.Range("D32").Value = "'0.178"
$ws.Range("E32").Value = "  +12.67%  "
$ws.Range("D33").Value = "'6.62"
$ws.Range("E33").Value = "  +7.96%  "
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("D35").Value = "'1.98"
$ws.Range("E35").Value = "  +6.08%  "
$ws.Range("D36").Value = "'26.54"
$ws.Range("E36").Value = "  +3.13%  "
$ws.Range("D37").Value = "3.616.82"
$ws.Range("E37").Value = "  +2.61%  "
$ws.Range("D38").Value = "'8.47"
$ws.Range("E38").Value = "  +5.61%  "
$ws.Range("D39").Value = "'2.44"
$ws.Range("E39").Value = "  +15.03%  "
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("D41").Value = "'0.0929"
$ws.Range("E41").Value = "  +8.18%  "
$ws.Range("D42").Value = "'0.998"
$ws.Range("E42").Value = "  -0.17%  "
$ws.Range("D43").Value = "'176.54"
$ws.Range("E43").Value = "  +1.42%  "
$ws.Range("D44").Value = "'5.64"
$ws.Range("E44").Value = "  +1.72%  "
$ws.Range("D45").Value = "'32.03"
$ws.Range("E45").Value = "  +17.98%  "
$ws.Range("D46").Value = "'0.913"
$ws.Range("E46").Value = "  +2.58%  "
$ws.Range("E47").Value = "  +12.53%  "
$ws.Range("D48").Value = "'2.83"
$ws.Range("E48").Value = "  +11.16%  "
$ws.Range("D49").Value = "'46.38"
$ws.Range("E49").Value = "  +2.26%  "
$ws.Range("D50").Value = "'7.82"
$ws.Range("E50").Value = "  +3.55%  "
$ws.Range("D51").Value = "'0.268"
$ws.Range("E51").Value = "  +7.61%  "
